$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 84, shifting existing rows 84:180 down to 85:181
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new data record
$ws.Cells.Item(84, 1).Value = 11
$ws.Cells.Item(84, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(84, 3).Value = "Bíobío"
$ws.Cells.Item(84, 4).Value = 44539
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 5).Value = 8
$ws.Cells.Item(84, 6).Value = 100112008
$ws.Cells.Item(84, 7).Value = "Coliflor"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 1800
$ws.Cells.Item(84, 11).Value = 550
$ws.Cells.Item(84, 12).Value = 600
$ws.Cells.Item(84, 13).Value = 572
$ws.Cells.Item(84, 14).Value = "`$/unidad"
$ws.Cells.Item(84, 15).Value = "Región Metropolitana"
$ws.Cells.Item(84, 16).Value = 572
$ws.Cells.Item(84, 17).Value = 1
$ws.Cells.Item(84, 18).Value = "Hortaliza"
